$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'20.532.28"
$ws.Range("E2").Value = '  +2.18%  '
$ws.Range("D3").Value = "'1.476.34"
$ws.Range("E3").Value = '  +3.77%  '
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = '  +0.75%  '
$ws.Range("D5").Value = "'0.9595"
$ws.Range("E5").Value = '  -3.73%  '
$ws.Range("D6").Value = "'276.79"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = "'0.3651"
$ws.Range("E7").Value = '  -1.45%  '
$ws.Range("D8").Value = "'0.3056"
$ws.Range("E8").Value = '  -2.82%  '
$ws.Range("D9").Value = "'39.74"
$ws.Range("E9").Value = '  +0.67%  '
$ws.Range("D10").Value = "'1.058"
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("D11").Value = "'0.06630"
$ws.Range("E11").Value = '  +1.26%  '
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").Value = "'18.21"
$ws.Range("E13").Value = '  +2.01%  '
$ws.Range("D14").Value = "'5.479"
$ws.Range("D15").Value = "'6.175"
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("D17").Value = "'1.475.30"
$ws.Range("E17").Value = '  +3.76%  '
$ws.Range("D18").Value = "'0.05906"
$ws.Range("E18").Value = '  +3.45%  '
$ws.Range("D19").Value = "'0.9650"
$ws.Range("E19").Value = '  -3.21%  '
$ws.Range("D20").Value = "'69.41"
$ws.Range("E20").Value = '  -3.28%  '
$ws.Range("D21").Value = "'5.472"
$ws.Range("E21").Value = '  -2.69%  '
$ws.Range("D22").Value = "'14.49"
$ws.Range("E22").Value = '  -2.66%  '
$ws.Range("D23").Value = "'11.01"
$ws.Range("E23").Value = '  -0.74%  '
$ws.Range("D24").Value = "'2.249"
$ws.Range("E24").Value = '  +1.12%  '
$ws.Range("D25").Value = "'20.593.24"
$ws.Range("E25").Value = '  +2.46%  '
$ws.Range("D26").Value = "'141.91"
$ws.Range("E26").Value = '  +5.73%  '
$ws.Range("D27").Value = "'2.136"
$ws.Range("E27").Value = '  -6.97%  '
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("D29").Value = "'1.632.02"
$ws.Range("D30").Value = "'113.81"
$ws.Range("E30").Value = '  +2.25%  '
$ws.Range("E31").Value = '  -1.18%  '
$ws.Range("D32").Value = "'0.8211"
$ws.Range("E32").Value = '  -0.64%  '
$ws.Range("D33").Value = "'4.981"
$ws.Range("E33").Value = '  -5.76%  '
$ws.Range("D34").Value = "'0.07946"
$ws.Range("E34").Value = '  +1.69%  '
$ws.Range("D35").Value = "'1.537"
$ws.Range("E35").Value = '  +4.14%  '
$ws.Range("D36").Value = "'1.237"
$ws.Range("E36").Value = '  +11.23%  '
$ws.Range("D37").Value = "'0.05761"
$ws.Range("E37").Value = '  -1.66%  '
$ws.Range("D38").Value = "'4.739"
$ws.Range("E38").Value = '  -3.81%  '
$ws.Range("E39").Value = '  -1.28%  '
$ws.Range("D40").Value = "'0.02039"
$ws.Range("E40").Value = '  -1.30%  '
$ws.Range("D41").Value = "'7.611"
$ws.Range("E41").Value = '  -5.26%  '
$ws.Range("E42").Value = '  -3.84%  '
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("D44").Value = "'0.5295"
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("D46").Value = "'12.12"
$ws.Range("E46").Value = '  -1.99%  '
$ws.Range("D47").Value = "'117.19"
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("D48").Value = "'0.5196"
$ws.Range("E48").Value = '  -0.82%  '
$ws.Range("D49").Value = "'1.782"
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").Value = "'0.06464"
$ws.Range("D51").Value = "'0.9921"
$ws.Range("E51").Value = '  -0.54%  '
